$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: sensor tests merged into a single "testSensor" case ---
$ws.Range("D12").Value = 'testSensor'
$ws.Range("E12").Value = 'Test adding a Sensor to a start trigger and one to a stop trigger. Test ability to toggle a Sensor on and off'
$ws.Range("H12").Value = 'Merged with previous tests startsensor and stopsensor'
$ws.Range("J12").Value = 'Associate a specific trigger to a sensor'

# --- Row 13: now holds the TestLane case (moved up from row 16) ---
$ws.Range("D13").Value = 'TestLane'
$ws.Range("E13").Value = 'Test Lane class functionality'
$ws.Range("H13").Clear()
$ws.Range("J13").Value = 'Handle multiple Racers at once, DNF cases, SWAP cases, '

# --- Row 14: now holds testRacerStartFinish (moved down from row 12) ---
$ws.Range("D14").Value = 'testRacerStartFinish'
$ws.Range("E14").Value = 'Test starting and stoping a Racer and associated getters '
$ws.Range("H14").Value = 'Merged with TestLane'
$ws.Range("J14").Value = 'Racer knows own time info'

# --- Row 15: new "testUpdateTime" test case, replacing TestToggleSensor ---
$ws.Range("D15").Value = 'testUpdateTime'
$ws.Range("E15").Value = 'Test ability to update time'
$ws.Range("G15").Clear()
$ws.Range("H15").Clear()
$ws.Range("J15").Clear()

# --- Row 16: new "TestNormalIND" case, replacing the old TestLane row ---
$ws.Range("D16").Value = 'TestNormalIND'
$ws.Range("E16").Value = 'Test creation and sample run of an individual event'
$ws.Range("F16").Clear()
$ws.Range("G16").Value = 'In progress'
$ws.Range("J16").Clear()

# --- Rows 17-19: brand new "TestNormal*" scenario rows ---
$ws.Range("D17").Value = 'TestNormalPARAIND'
$ws.Range("E17").Value = 'Test creation and sample run of a paraividual event'
$ws.Range("G17").Value = 'In progress'

$ws.Range("D18").Value = 'TestNormalGRP'
$ws.Range("E18").Value = 'Test a group event'
$ws.Range("G18").Value = 'In progress'

$ws.Range("D19").Value = 'TestNormalPARAGRP'
$ws.Range("E19").Value = 'Test a paragrpividual event'
$ws.Range("G19").Value = 'In progress'

# Row heights that Excel recorded as explicit ("custom") after the content grew
$ws.Rows.Item(12).RowHeight = 48
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 15.75

# Reset the view: scroll back to the top-left and select G8 (matches the saved view)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G8").Select()

# Workbook-level metadata touched by the resave
$wb.Windows.Item(1).Height = 9090
